$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": Latest HO Xliff Generate Date for 4a35a5b2-... and a406add0-...
#     both rows previously showed 2016-08-18 06:15:01, now show 2016-08-18 06:15:57
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-18 06:15:57"
$wsOverview.Range("G4").Value = "2016-08-18 06:15:57"

# --- Sheet "zh-cn": Priority ht -> mt, and updated handoff/handback datetimes
#     for the 4a35a5b2-... and a406add0-... rows (rows 3 and 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-18 06:15:52"
$wsZhCn.Range("H4").Value = "2016-08-18 06:15:52"
$wsZhCn.Range("K3").Value = "2016-08-18 06:16:14"
$wsZhCn.Range("K4").Value = "2016-08-18 06:16:14"

# --- Sheet "de-de": Priority ht -> mt (shares string with zh-cn sheet), and
#     updated Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-18 06:15:57"
$wsDeDe.Range("H4").Value = "2016-08-18 06:15:57"
$wsDeDe.Range("K3").Value = "2016-08-18 06:16:21"
$wsDeDe.Range("K4").Value = "2016-08-18 06:16:21"
